{"js": "// Add a new \"Version updated to 1.1.0\" bullet item right after the\n// existing \"Version updated to 1.0.14\" item at the end of the document,\n// matching the same list formatting (ListParagraph style, numId 6).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document is \"Version updated to 1.0.14\".\nconst lastPara = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the new paragraph right after it; it inherits the ListParagraph\n// style + list numbering from the paragraph it's split/inserted from.\nlastPara.insertParagraph(\"Version updated to 1.1.0\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new \"Version updated to 1.1.0\" bullet item right after the\n# existing \"Version updated to 1.0.14\" item at the end of the document,\n# matching the same list formatting (ListParagraph style, numId 6).\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is \"Version updated to 1.0.14\".\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertParagraphAfter()\n\n# Insert the new text into the freshly created paragraph (it inherited\n# the ListParagraph style + list numbering from the split paragraph).\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Version updated to 1.1.0\"\n"}
